# apply phero system in ShieldAbility
# Reduce the MaxHp (column J) value for the first enemy row from 100 to 70,
# and move the active selection to M9 to match the author's last cursor
# position when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MaxHp value for row 2 (J2): 100 -> 70
$ws.Range("J2").Value = 70

# Move selection to M9 (matches saved cursor position in the diff)
$ws.Range("M9").Select()
